# ---------------------------------------------------------------------------
# Requirements_second_turn_Aseev_Malofeeva_v1.1.docx -- "update requirments
# and resume"
#
#   1. Drop the stray "_GoBack" bookmark that sits right after the item
#      ending "...решены) пользователем." (an artifact of the cursor
#      position at the time the document was last saved).
#   2. Add a brand-new requirement bullet right after the one ending
#      "...хранимое в формате XML." -- "Объем загружаемого пользователем
#      файла ... не должен превышать 64Кб." -- and this new bullet now
#      carries a fresh "_GoBack" bookmark (the cursor position when this
#      revision was saved).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# wdGoToBookmark-ish direct lookup: Bookmarks.Item("_GoBack") resolves to
# the correct (hidden, auto-generated) bookmark object positionally even
# though its own .Name getter is unreliable for this particular bookmark,
# so .Delete() on it removes exactly the right <w:bookmarkStart/End> pair.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Locate the end of "...хранимое в формате XML." and split a new paragraph
# off right after it. InsertParagraphAfter() clones the paragraph/run
# formatting from the split point, which is exactly the a3 / numId=14 /
# ilvl=1 / Times New Roman list formatting the surrounding bullets use.
$anchor = $d.Content
$null = $anchor.Find.Execute("хранимое в формате XML.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Collapse(0)   # wdCollapseEnd
$anchor.InsertParagraphAfter()

# The freshly split paragraph is now the last paragraph in the document
# (it was inserted right before the trailing empty closing paragraph).
$newParaIndex = $d.Paragraphs.Count - 1
$newPara = $d.Paragraphs.Item($newParaIndex)

# Fill in the new bullet's text (leave the paragraph mark alone) -- append
# a one-off placeholder character so the bookmark-anchoring range below is
# never collapsed-at-paragraph-end (a configuration Bookmarks.Add mishandles).
$body = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$body.Text = " Объем загружаемого пользователем файла с исходным кодом для решения задачи не должен превышать 64Кб.X"

# Anchor a new "_GoBack" bookmark at the very end of the new bullet's text.
# Bookmarks.Add on a truly collapsed range sitting exactly at a
# paragraph-end boundary drops the bookmark at document position 0 instead
# (an engine quirk), so bracket the placeholder character with a
# non-collapsed range -- that anchors correctly -- then erase just the
# placeholder character's text; the bookmark tags stay right where they
# were placed, now collapsed at the true end of the paragraph's text.
$placeholder = $d.Range($body.End - 1, $body.End)
$d.Bookmarks.Add("_GoBack", $placeholder)
$placeholder2 = $d.Range($body.End - 1, $body.End)
$placeholder2.Text = ""
